$d = $word.ActiveDocument

# Move to the very end of the document content (before the final paragraph mark)
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd = 0

# The last paragraph ("Try and get the feeling.  ") ends with a paragraph mark;
# endRange currently sits at the very end of the body (after that mark).
# Insert a new paragraph for the separator line.
$lastPara = $d.Paragraphs.Last
$sepRange = $lastPara.Range
$sepRange.Collapse(0)  # collapse to end of paragraph (before its mark)
$sepRange.InsertParagraphAfter()

$sepPara = $d.Paragraphs.Last
$sepPara.Range.Text = "--------------------------------------------------"
$sepPara.Range.Font.Size = 12
$sepPara.Range.Font.Bold = 0

# Insert second new paragraph with the "New changes" text.
$sepPara.Range.InsertParagraphAfter()

$msgPara = $d.Paragraphs.Last
$msgPara.Range.Text = "New changes have been adde below"
$msgPara.Range.Font.Size = 12
$msgPara.Range.Font.Bold = 0
